$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1157.2903
$ws.Range("I33").Value = 961.14813
$ws.Range("J33").Value = 2481.25
$ws.Range("K33").Value = 961.14813
$ws.Range("L33").Value = 2481.25
$ws.Range("M33").Value = -732.14813
$ws.Range("N33").Value = -2939.25

$ws.Range("H70").Value = 1399.6
$ws.Range("I70").Value = 1399.6
$ws.Range("K70").Value = 4198.799999999999
$ws.Range("M70").Value = -3928.799999999999

$ws.Range("H73").Value = 1399.6
$ws.Range("I73").Value = 1399.6
$ws.Range("K73").Value = 4198.799999999999
$ws.Range("M73").Value = -3262.799999999999

$ws.Range("H137").Value = 5487.778
$ws.Range("I137").Value = 2271.6365
$ws.Range("J137").Value = 10541.714
$ws.Range("K137").Value = 6814.9095
$ws.Range("L137").Value = 31625.142
$ws.Range("M137").Value = -4264.9095
$ws.Range("N137").Value = -36725.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 19276136
$ws.Range("I61").Value = 26321274
$ws.Range("J61").Value = 153619.72
$ws.Range("K61").Value = 26321274
$ws.Range("L61").Value = 153619.72
$ws.Range("M61").Value = -26321062
$ws.Range("N61").Value = -154043.72

$ws.Range("H74").Value = 13901483
$ws.Range("I74").Value = 25002664
$ws.Range("J74").Value = 25008.375
$ws.Range("K74").Value = 25002664
$ws.Range("L74").Value = 25008.375
$ws.Range("M74").Value = -25001790
$ws.Range("N74").Value = -26756.375

$ws.Range("H77").Value = 13901483
$ws.Range("I77").Value = 25002664
$ws.Range("J77").Value = 25008.375
$ws.Range("K77").Value = 125013320
$ws.Range("L77").Value = 125041.875
$ws.Range("M77").Value = -125008952
$ws.Range("N77").Value = -133777.875

$ws.Range("H105").Value = 49750
$ws.Range("J105").Value = 49750
$ws.Range("L105").Value = 49750
$ws.Range("N105").Value = -56738

$ws.Range("H122").Value = 1699.8948
$ws.Range("I122").Value = 1138.3846
$ws.Range("K122").Value = 3415.1538
$ws.Range("M122").Value = -965.1538

$ws.Range("H132").Value = 7917.4
$ws.Range("I132").Value = 5677.3335
$ws.Range("J132").Value = 13377.5625
$ws.Range("K132").Value = 17032.0005
$ws.Range("L132").Value = 40132.6875
$ws.Range("M132").Value = -14502.0005
$ws.Range("N132").Value = -45192.6875

$ws.Range("H133").Value = 79997.25
$ws.Range("J133").Value = 79997.25
$ws.Range("L133").Value = 79997.25
$ws.Range("N133").Value = -85057.25

$ws.Range("H136").Value = 19276136
$ws.Range("I136").Value = 26321274
$ws.Range("J136").Value = 153619.72
$ws.Range("K136").Value = 78963822
$ws.Range("L136").Value = 460859.16
$ws.Range("M136").Value = -78961272
$ws.Range("N136").Value = -465959.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 154.25
$ws.Range("I22").Value = 153.42857
$ws.Range("K22").Value = 153.42857
$ws.Range("M22").Value = 19.57142999999999

$ws.Range("H86").Value = 2825.1
$ws.Range("I86").Value = 2836.4285
$ws.Range("K86").Value = 2836.4285
$ws.Range("M86").Value = -1713.4285

$ws.Range("H89").Value = 2825.1
$ws.Range("I89").Value = 2836.4285
$ws.Range("K89").Value = 14182.1425
$ws.Range("M89").Value = -8566.1425

$ws.Range("H96").Value = 32214.3
$ws.Range("J96").Value = 61411
$ws.Range("L96").Value = 61411
$ws.Range("N96").Value = -66903

$ws.Range("H112").Value = 130000
$ws.Range("J112").Value = 130000
$ws.Range("L112").Value = 130000
$ws.Range("N112").Value = -132954

$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3373.75
$ws.Range("I62").Value = 3995
$ws.Range("J62").Value = 3166.6667
$ws.Range("K62").Value = 3995
$ws.Range("L62").Value = 3166.6667
$ws.Range("M62").Value = -3371
$ws.Range("N62").Value = -4414.6667

$ws.Range("H65").Value = 3373.75
$ws.Range("I65").Value = 3995
$ws.Range("J65").Value = 3166.6667
$ws.Range("K65").Value = 19975
$ws.Range("L65").Value = 15833.3335
$ws.Range("M65").Value = -16855
$ws.Range("N65").Value = -22073.3335

$ws.Range("H134").Value = 316180.2
$ws.Range("I134").Value = 418562.47
$ws.Range("K134").Value = 1255687.41
$ws.Range("M134").Value = -1253152.41

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9151.3125
$ws.Range("I3").Value = 3269.7144
$ws.Range("J3").Value = 13725.889
$ws.Range("K3").Value = 9809.143199999999
$ws.Range("L3").Value = 41177.667
$ws.Range("M3").Value = -9697.143199999999
$ws.Range("N3").Value = -41401.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6313.4
$ws.Range("I102").Value = 4305
$ws.Range("K102").Value = 4305
$ws.Range("M102").Value = -2683

$ws.Range("H132").Value = 33336992
$ws.Range("I132").Value = 37040750
$ws.Range("K132").Value = 111122250
$ws.Range("M132").Value = -111119720

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 950.1429
$ws.Range("I82").Value = 374
$ws.Range("J82").Value = 1180.6
$ws.Range("K82").Value = 374
$ws.Range("L82").Value = 1180.6
$ws.Range("M82").Value = -13
$ws.Range("N82").Value = -1902.6

$ws.Range("H85").Value = 950.1429
$ws.Range("I85").Value = 374
$ws.Range("J85").Value = 1180.6
$ws.Range("K85").Value = 374
$ws.Range("L85").Value = 1180.6
$ws.Range("M85").Value = 874
$ws.Range("N85").Value = -3676.6

$ws.Range("I93").Value = 66668396
$ws.Range("K93").Value = 66668396
$ws.Range("M93").Value = -66667148

$ws.Range("H104").Value = 43748.2
$ws.Range("J104").Value = 43748.2
$ws.Range("L104").Value = 43748.2
$ws.Range("N104").Value = -50736.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 6727.1816
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -710

$ws.Range("H98").Value = 67542.5
$ws.Range("J98").Value = 67542.5
$ws.Range("L98").Value = 67542.5
$ws.Range("N98").Value = -73532.5

$ws.Range("H132").Value = 12052.777
$ws.Range("I132").Value = 1210.7142
$ws.Range("K132").Value = 3632.1426
$ws.Range("M132").Value = -1102.1426
